# Scheduled data refresh: update market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on each Leve sheet
# to the latest pulled values. Rows where HQ price/profit no longer
# applies have those cells cleared instead of zeroed.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 364.54544
$ws.Range("I5").Value = 363.66666
$ws.Range("J5").Value = 365.6
$ws.Range("K5").Value = 363.66666
$ws.Range("L5").Value = 365.6
$ws.Range("M5").Value = -248.66666
$ws.Range("N5").Value = -595.6

$ws.Range("H20").Value = 12040
$ws.Range("I20").Value = 7066.6665
$ws.Range("J20").Value = 19500
$ws.Range("K20").Value = 7066.6665
$ws.Range("L20").Value = 19500
$ws.Range("M20").Value = -6836.6665
$ws.Range("N20").Value = -19960

$ws.Range("H35").Value = 12040
$ws.Range("I35").Value = 7066.6665
$ws.Range("J35").Value = 19500
$ws.Range("K35").Value = 7066.6665
$ws.Range("L35").Value = 19500
$ws.Range("M35").Value = -6687.6665
$ws.Range("N35").Value = -20258

$ws.Range("H87").Value = 99999.766
$ws.Range("J87").Value = 99999.766
$ws.Range("L87").Value = 99999.766
$ws.Range("N87").Value = -102495.766

$ws.Range("H90").Value = 99999.766
$ws.Range("J90").Value = 99999.766
$ws.Range("L90").Value = 299999.298
$ws.Range("N90").Value = -312479.298

$ws.Range("H92").Value = 2585.3333
$ws.Range("I92").Value = 2377
$ws.Range("J92").Value = 3002
$ws.Range("K92").Value = 2377
$ws.Range("L92").Value = 3002
$ws.Range("M92").Value = -1129
$ws.Range("N92").Value = -5498

$ws.Range("H100").Value = 6623.3965
$ws.Range("I100").Value = 1315.3334
$ws.Range("K100").Value = 1315.3334
$ws.Range("M100").Value = -774.3334

$ws.Range("H112").Value = 2520.25
$ws.Range("J112").Value = 2520.25
$ws.Range("L112").Value = 7560.75
$ws.Range("N112").Value = -9776.75

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()

$ws.Range("H138").Value = 3586.8774
$ws.Range("I138").Value = 1448.4166
$ws.Range("J138").Value = 5639.8
$ws.Range("K138").Value = 4345.2498
$ws.Range("L138").Value = 16919.4
$ws.Range("M138").Value = 794.7502000000004
$ws.Range("N138").Value = -27199.4


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6135805
$ws.Range("I2").Value = 8365343.5
$ws.Range("K2").Value = 8365343.5
$ws.Range("M2").Value = -8365230.5

$ws.Range("H31").Value = 13333
$ws.Range("I31").Value = 13333
$ws.Range("K31").Value = 13333
$ws.Range("M31").Value = -13039

$ws.Range("H97").Value = 2646910.8
$ws.Range("I97").Value = 3368545.5
$ws.Range("J97").Value = 916.3333
$ws.Range("K97").Value = 3368545.5
$ws.Range("L97").Value = 916.3333
$ws.Range("M97").Value = -3368049.5
$ws.Range("N97").Value = -1908.3333

$ws.Range("H116").Value = 6135805
$ws.Range("I116").Value = 8365343.5
$ws.Range("K116").Value = 8365343.5
$ws.Range("M116").Value = -8363049.5

$ws.Range("H122").Value = 2733.138
$ws.Range("I122").Value = 2842.5925
$ws.Range("K122").Value = 8527.7775
$ws.Range("M122").Value = -6077.7775


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6135805
$ws.Range("I3").Value = 8365343.5
$ws.Range("K3").Value = 8365343.5
$ws.Range("M3").Value = -8365229.5

$ws.Range("H20").Value = 2598.2104
$ws.Range("I20").Value = 3546.111
$ws.Range("K20").Value = 3546.111
$ws.Range("M20").Value = -3299.111

$ws.Range("H107").Value = 1218.5428
$ws.Range("I107").Value = 1071.48
$ws.Range("K107").Value = 1071.48
$ws.Range("M107").Value = 848.52


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 333434.66
$ws.Range("I23").Value = 333434.66
$ws.Range("K23").Value = 1000303.98
$ws.Range("M23").Value = -1000068.98

$ws.Range("H113").Value = 1636.16
$ws.Range("I113").Value = 1121
$ws.Range("J113").Value = 2111.6924
$ws.Range("K113").Value = 3363
$ws.Range("L113").Value = 6335.0772
$ws.Range("M113").Value = -1193
$ws.Range("N113").Value = -10675.0772

$ws.Range("H132").Value = 4343.4443
$ws.Range("I132").Value = 8399
$ws.Range("J132").Value = 3184.7144
$ws.Range("K132").Value = 75591
$ws.Range("L132").Value = 28662.4296
$ws.Range("M132").Value = -73061
$ws.Range("N132").Value = -33722.4296


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1125
$ws.Range("I97").Value = 930.8182
$ws.Range("K97").Value = 930.8182
$ws.Range("M97").Value = -434.8182

$ws.Range("H102").Value = 4374.091
$ws.Range("I102").Value = 4510.1
$ws.Range("K102").Value = 4510.1
$ws.Range("M102").Value = -2888.1

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3365.2683
$ws.Range("I22").Value = 2076.7144
$ws.Range("J22").Value = 4718.25
$ws.Range("K22").Value = 2076.7144
$ws.Range("L22").Value = 4718.25
$ws.Range("M22").Value = -1781.7144
$ws.Range("N22").Value = -5308.25

$ws.Range("H27").Value = 3365.2683
$ws.Range("I27").Value = 2076.7144
$ws.Range("J27").Value = 4718.25
$ws.Range("K27").Value = 2076.7144
$ws.Range("L27").Value = 4718.25
$ws.Range("M27").Value = -1969.7144
$ws.Range("N27").Value = -4932.25

$ws.Range("H68").Value = 3739.5625
$ws.Range("I68").Value = 2124.25
$ws.Range("J68").Value = 5354.875
$ws.Range("K68").Value = 2124.25
$ws.Range("L68").Value = 5354.875
$ws.Range("M68").Value = -1375.25
$ws.Range("N68").Value = -6852.875

$ws.Range("H71").Value = 3739.5625
$ws.Range("I71").Value = 2124.25
$ws.Range("J71").Value = 5354.875
$ws.Range("K71").Value = 10621.25
$ws.Range("L71").Value = 26774.375
$ws.Range("M71").Value = -6877.25
$ws.Range("N71").Value = -34262.375

$ws.Range("H93").Value = 1187.7273
$ws.Range("I93").Value = 1187.7273
$ws.Range("K93").Value = 1187.7273
$ws.Range("M93").Value = 60.27269999999999

$ws.Range("H122").Value = 15629802
$ws.Range("I122").Value = 17246302
$ws.Range("J122").Value = 3633.3333
$ws.Range("K122").Value = 51738906
$ws.Range("L122").Value = 10899.9999
$ws.Range("M122").Value = -51736456
$ws.Range("N122").Value = -15799.9999

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1894.4286
$ws.Range("I122").Value = 1894.4286
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5683.2858
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3233.2858
$ws.Range("N122").ClearContents()

